$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3333.25
$ws.Cells.Item(64, 9).Value = 3304.7273
$ws.Cells.Item(64, 10).Value = 3396
$ws.Cells.Item(64, 11).Value = 3304.7273
$ws.Cells.Item(64, 12).Value = 3396
$ws.Cells.Item(64, 13).Value = -3056.7273
$ws.Cells.Item(64, 14).Value = -3892

$ws.Cells.Item(67, 8).Value = 3333.25
$ws.Cells.Item(67, 9).Value = 3304.7273
$ws.Cells.Item(67, 10).Value = 3396
$ws.Cells.Item(67, 11).Value = 3304.7273
$ws.Cells.Item(67, 12).Value = 3396
$ws.Cells.Item(67, 13).Value = -2446.7273
$ws.Cells.Item(67, 14).Value = -5112

$ws.Cells.Item(132, 8).Value = 1205.5962
$ws.Cells.Item(132, 9).Value = 1145.5918
$ws.Cells.Item(132, 11).Value = 3436.7754
$ws.Cells.Item(132, 13).Value = -906.7753999999995

$ws.Cells.Item(134, 8).Value = 45000
$ws.Cells.Item(134, 10).Value = 45000
$ws.Cells.Item(134, 12).Value = 45000
$ws.Cells.Item(134, 14).Value = -55140

$ws.Cells.Item(140, 8).Value = 71947.62
$ws.Cells.Item(140, 10).Value = 92726.664
$ws.Cells.Item(140, 12).Value = 92726.664
$ws.Cells.Item(140, 14).Value = -103086.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2618.3572
$ws.Cells.Item(61, 9).Value = 2525.3447
$ws.Cells.Item(61, 10).Value = 2825.8462
$ws.Cells.Item(61, 11).Value = 2525.3447
$ws.Cells.Item(61, 12).Value = 2825.8462
$ws.Cells.Item(61, 13).Value = -2313.3447
$ws.Cells.Item(61, 14).Value = -3249.8462

$ws.Cells.Item(74, 8).Value = 686.7727
$ws.Cells.Item(74, 9).Value = 505.73685
$ws.Cells.Item(74, 10).Value = 1833.3334
$ws.Cells.Item(74, 11).Value = 505.73685
$ws.Cells.Item(74, 12).Value = 1833.3334
$ws.Cells.Item(74, 13).Value = 368.26315
$ws.Cells.Item(74, 14).Value = -3581.3334

$ws.Cells.Item(77, 8).Value = 686.7727
$ws.Cells.Item(77, 9).Value = 505.73685
$ws.Cells.Item(77, 10).Value = 1833.3334
$ws.Cells.Item(77, 11).Value = 2528.68425
$ws.Cells.Item(77, 12).Value = 9166.666999999999
$ws.Cells.Item(77, 13).Value = 1839.31575
$ws.Cells.Item(77, 14).Value = -17902.667

$ws.Cells.Item(136, 8).Value = 2618.3572
$ws.Cells.Item(136, 9).Value = 2525.3447
$ws.Cells.Item(136, 10).Value = 2825.8462
$ws.Cells.Item(136, 11).Value = 7576.034100000001
$ws.Cells.Item(136, 12).Value = 8477.5386
$ws.Cells.Item(136, 13).Value = -5026.034100000001
$ws.Cells.Item(136, 14).Value = -13577.5386

$ws.Cells.Item(138, 8).Value = 67712.5
$ws.Cells.Item(138, 10).Value = 67712.5
$ws.Cells.Item(138, 12).Value = 67712.5
$ws.Cells.Item(138, 14).Value = -77992.5

$ws.Cells.Item(139, 8).Value = 54722.5
$ws.Cells.Item(139, 10).Value = 54722.5
$ws.Cells.Item(139, 12).Value = 54722.5
$ws.Cells.Item(139, 14).Value = -65002.5

$ws.Cells.Item(140, 8).Value = 105000
$ws.Cells.Item(140, 10).Value = 105000
$ws.Cells.Item(140, 12).Value = 105000
$ws.Cells.Item(140, 14).Value = -115360

$ws.Cells.Item(141, 8).Value = 60992.855
$ws.Cells.Item(141, 10).Value = 63992.31
$ws.Cells.Item(141, 12).Value = 63992.31
$ws.Cells.Item(141, 14).Value = -74352.31

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 88900
$ws.Cells.Item(140, 10).Value = 88900
$ws.Cells.Item(140, 12).Value = 88900
$ws.Cells.Item(140, 14).Value = -99260

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10578.647
$ws.Cells.Item(31, 9).Value = 3798.2856
$ws.Cells.Item(31, 11).Value = 3798.2856
$ws.Cells.Item(31, 13).Value = -3503.2856

$ws.Cells.Item(34, 8).Value = 10578.647
$ws.Cells.Item(34, 9).Value = 3798.2856
$ws.Cells.Item(34, 11).Value = 3798.2856
$ws.Cells.Item(34, 13).Value = -3596.2856

$ws.Cells.Item(62, 8).Value = 4041.2632
$ws.Cells.Item(62, 9).Value = 4052.2666
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 11).Value = 4052.2666
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 13).Value = -3428.2666
$ws.Cells.Item(62, 14).Value = -5248

$ws.Cells.Item(65, 8).Value = 4041.2632
$ws.Cells.Item(65, 9).Value = 4052.2666
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 11).Value = 20261.333
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 13).Value = -17141.333
$ws.Cells.Item(65, 14).Value = -26240

$ws.Cells.Item(86, 8).Value = 28228.928
$ws.Cells.Item(86, 9).Value = 3098.9524
$ws.Cells.Item(86, 10).Value = 103618.86
$ws.Cells.Item(86, 11).Value = 3098.9524
$ws.Cells.Item(86, 12).Value = 103618.86
$ws.Cells.Item(86, 13).Value = -1975.9524
$ws.Cells.Item(86, 14).Value = -105864.86

$ws.Cells.Item(89, 8).Value = 28228.928
$ws.Cells.Item(89, 9).Value = 3098.9524
$ws.Cells.Item(89, 10).Value = 103618.86
$ws.Cells.Item(89, 11).Value = 15494.762
$ws.Cells.Item(89, 12).Value = 518094.3
$ws.Cells.Item(89, 13).Value = -9878.762000000001
$ws.Cells.Item(89, 14).Value = -529326.3

$ws.Cells.Item(138, 8).Value = 49275
$ws.Cells.Item(138, 10).Value = 49275
$ws.Cells.Item(138, 12).Value = 49275
$ws.Cells.Item(138, 14).Value = -59555

$ws.Cells.Item(140, 8).Value = 64000
$ws.Cells.Item(140, 10).Value = 64000
$ws.Cells.Item(140, 12).Value = 64000
$ws.Cells.Item(140, 14).Value = -74360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 5209157
$ws.Cells.Item(131, 9).Value = 829.9167
$ws.Cells.Item(131, 10).Value = 6945266
$ws.Cells.Item(131, 11).Value = 2489.7501
$ws.Cells.Item(131, 12).Value = 20835798
$ws.Cells.Item(131, 13).Value = 2550.2499
$ws.Cells.Item(131, 14).Value = -20845878

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 69033.336
$ws.Cells.Item(138, 10).Value = 69033.336
$ws.Cells.Item(138, 12).Value = 69033.336
$ws.Cells.Item(138, 14).Value = -79313.336

$ws.Cells.Item(140, 8).Value = 89989
$ws.Cells.Item(140, 10).Value = 89989
$ws.Cells.Item(140, 12).Value = 89989
$ws.Cells.Item(140, 14).Value = -100349

$ws.Cells.Item(141, 8).Value = 43666.668
$ws.Cells.Item(141, 10).Value = 43666.668
$ws.Cells.Item(141, 12).Value = 43666.668
$ws.Cells.Item(141, 14).Value = -54026.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2635
$ws.Cells.Item(136, 9).Value = 2534.0952
$ws.Cells.Item(136, 10).Value = 2846.9
$ws.Cells.Item(136, 11).Value = 7602.285600000001
$ws.Cells.Item(136, 12).Value = 8540.700000000001
$ws.Cells.Item(136, 13).Value = -5052.285600000001
$ws.Cells.Item(136, 14).Value = -13640.7

$ws.Cells.Item(138, 8).Value = 56007.168
$ws.Cells.Item(138, 10).Value = 56007.168
$ws.Cells.Item(138, 12).Value = 56007.168
$ws.Cells.Item(138, 14).Value = -66287.16800000001

$ws.Cells.Item(139, 8).Value = 53216.668
$ws.Cells.Item(139, 10).Value = 62860
$ws.Cells.Item(139, 12).Value = 62860
$ws.Cells.Item(139, 14).Value = -73140

$ws.Cells.Item(140, 8).Value = 58480
$ws.Cells.Item(140, 10).Value = 59000
$ws.Cells.Item(140, 12).Value = 59000
$ws.Cells.Item(140, 14).Value = -69360

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 80000
$ws.Cells.Item(139, 10).Value = 80000
$ws.Cells.Item(139, 12).Value = 80000
$ws.Cells.Item(139, 14).Value = -90280

$ws.Cells.Item(140, 8).Value = 29900
$ws.Cells.Item(140, 10).Value = 29900
$ws.Cells.Item(140, 12).Value = 29900
$ws.Cells.Item(140, 14).Value = -40260

$ws.Cells.Item(141, 8).Value = 73943
$ws.Cells.Item(141, 10).Value = 73943
$ws.Cells.Item(141, 12).Value = 73943
$ws.Cells.Item(141, 14).Value = -84303
